# Apply the Thu Oct  3 09:27:07 UTC 2024 "cryptos list" refresh.
# Updates Price (D) / Volume(1h) (E) for every ranked coin, and for the
# four rows whose rank order swapped (26/27, 29/30, 37/39) also rewrites
# the Coin name (B) and Link (C) so the row carries the new coin.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.626.84"
$ws.Range("E2").Value = "  -1.57%  "

# Row 3
$ws.Range("D3").Value = "2.345.16"
$ws.Range("E3").Value = "  -5.04%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "542.60"
$ws.Range("E5").Value = "  -1.96%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.34"
$ws.Range("E6").Value = "  -6.74%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.16%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.517"
$ws.Range("E8").Value = "  -13.32%  "

# Row 9
$ws.Range("D9").Value = "2.345.18"
$ws.Range("E9").Value = "  -4.86%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.103"
$ws.Range("E10").Value = "  -4.59%  "

# Row 11
$ws.Range("E11").Value = "  -0.20%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.18"
$ws.Range("E12").Value = "  -5.30%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.337"
$ws.Range("E13").Value = "  -5.33%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.48"
$ws.Range("E14").Value = "  -7.11%  "

# Row 15
$ws.Range("D15").Value = "2.760.81"
$ws.Range("E15").Value = "  -5.37%  "

# Row 16
$ws.Range("D16").Value = "60.303.25"
$ws.Range("E16").Value = "  -1.97%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000159"
$ws.Range("E17").Value = "  -4.82%  "

# Row 18
$ws.Range("D18").Value = "2.336.66"
$ws.Range("E18").Value = "  -5.45%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.49"
$ws.Range("E19").Value = "  -6.11%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.05"
$ws.Range("E20").Value = "  -4.15%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "311.70"
$ws.Range("E21").Value = "  -2.98%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.50"
$ws.Range("E22").Value = "  -9.88%  "

# Row 23
$ws.Range("E23").Value = "  -0.29%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.87"
$ws.Range("E24").Value = "  -2.10%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "62.57"
$ws.Range("E25").Value = "  -2.80%  "

# Row 26
$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.06"
$ws.Range("E26").Value = "  +3.43%  "

# Row 27
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.19%  "

# Row 28
$ws.Range("D28").Value = "2.449.30"
$ws.Range("E28").Value = "  -6.12%  "

# Row 29
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.90"
$ws.Range("E29").Value = "  -5.26%  "

# Row 30
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0882"
$ws.Range("E30").Value = "  -11.62%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "502.17"
$ws.Range("E31").Value = "  -10.97%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.38"
$ws.Range("E32").Value = "  -8.54%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.143"
$ws.Range("E33").Value = "  -4.24%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.78"
$ws.Range("E34").Value = "  -7.19%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.53"
$ws.Range("E35").Value = "  -4.88%  "

# Row 36
$ws.Range("E36").Value = "  -0.13%  "

# Row 37
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.40"
$ws.Range("E37").Value = "  -0.60%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.50"
$ws.Range("E38").Value = "  -8.14%  "

# Row 39
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.370"
$ws.Range("E39").Value = "  -3.02%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.23"
$ws.Range("E40").Value = "  -12.11%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.78"
$ws.Range("E41").Value = "  +0.88%  "

# Row 42
$ws.Range("E42").Value = "  +0.00%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "137.96"
$ws.Range("E43").Value = "  -4.45%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.12"
$ws.Range("E44").Value = "  -1.18%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "137.58"
$ws.Range("E45").Value = "  -6.87%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.53"
$ws.Range("E46").Value = "  -2.99%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.06"
$ws.Range("E47").Value = "  -15.69%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0507"
$ws.Range("E48").Value = "  -6.13%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.30"
$ws.Range("E49").Value = "  -11.15%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.564"
$ws.Range("E50").Value = "  -5.25%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0891"
$ws.Range("E51").Value = "  -5.21%  "
